$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("floor_2")

# --- Row 43: new 1x8 board entry ---
$ws.Range("L43").Value = 120
$ws.Range("M43").Formula = '=INT(L43/12) & " ft "'
$ws.Range("N43").Value = 1
$ws.Range("O43").Value = 6
$ws.Range("P43").Formula = "=O43*N43"
$ws.Range("Q43").Value = 6
$ws.Range("R43").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 44 ---
$ws.Range("L44").Value = 96
$ws.Range("M44").Formula = '=INT(L44/12) & " ft "'
$ws.Range("N44").Value = 5
$ws.Range("O44").Value = 4.8
$ws.Range("P44").Formula = "=O44*N44"
$ws.Range("R44").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 45 ---
$ws.Range("L45").Value = 144
$ws.Range("M45").Formula = '=INT(L45/12) & " ft "'
$ws.Range("N45").Value = 1
$ws.Range("O45").Value = 7.6
$ws.Range("P45").Formula = "=O45*N45"
$ws.Range("R45").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 46 ---
$ws.Range("L46").Value = 168
$ws.Range("M46").Formula = '=INT(L46/12) & " ft "'
$ws.Range("N46").Value = 2
$ws.Range("O46").Value = 8.86
$ws.Range("P46").Formula = "=O46*N46"
$ws.Range("R46").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 47: subtotal ---
$ws.Range("Q47").Formula = "=SUM(P44:P46)"



# --- Sheet view state ---
$ws.Range("P69").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 7
